$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Defined names: remove unused/obsolete ones (check, deli, deliberacao12) ---
foreach ($nm in @("check", "deli", "deliberacao12")) {
    foreach ($n in @($wb.Names)) {
        if ($n.Name -eq $nm) {
            $n.Delete()
        }
    }
}

# --- Text updates ---
$ws.Range("B11").Value = "Atualização do Plano de Desenvolvimento Urbano Integrado da Região Metropolitana de Belo Horizonte – PDUI-RMBH"
$ws.Range("B43").Value = "Contratações temporárias e outras despesas de pessoal"

# --- Value updates ---
$ws.Range("D11").Value = 4400000
$ws.Range("D13").Value = 36712000
$ws.Range("D28").Value = 40729352.11
$ws.Range("D51").Value = 2427295557.9000001

# --- New column widths / formatting touches (columns F and I get used) ---
$ws.Columns("F").ColumnWidth = 16.5703125
$ws.Columns("I").ColumnWidth = 25.85546875

# New blank formatted cells introduced by the edit
$ws.Range("I33").NumberFormat = "0.000"
$ws.Range("I33").WrapText = $true

$ws.Range("F59").NumberFormat = "_-* #,##0.00_-;-* #,##0.00_-;_-* ""-""??_-;_-@_-"
$ws.Range("F59").WrapText = $true

# Row height tweaks to match the author's re-save of rows 60-61
$ws.Rows(60).RowHeight = 30.75
$ws.Rows(61).RowHeight = 30.75

# Update the active selection to match the author's cursor position on save
$ws.Range("E6").Select()
